$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the style from the cell above so the new date cell matches existing formatting
$ws.Range("A37").Copy()
$ws.Range("A38").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add new row of data: date, hours, comment
$ws.Range("A38").Value2 = 40274
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = "Weekly Meeting"

# Update active selection to reflect new last row
$ws.Range("A39").Select()
